$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.963.59"
$ws.Range("E2").Value = "  -2.82%  "
$ws.Range("D3").Value = "1.888.59"
$ws.Range("E3").Value = "  -3.64%  "
$ws.Range("E4").Value = "  -0.89%  "
$ws.Range("D5").Value = "325.86"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("D7").Value = "0.4571"
$ws.Range("E7").Value = "  -4.10%  "
$ws.Range("D8").Value = "0.3935"
$ws.Range("E8").Value = "  -2.40%  "
$ws.Range("D9").Value = "50.40"
$ws.Range("E9").Value = "  -6.59%  "
$ws.Range("D10").Value = "0.08201"
$ws.Range("E10").Value = "  -3.28%  "
$ws.Range("D11").Value = "1.036"
$ws.Range("E11").Value = "  -2.37%  "
$ws.Range("D12").Value = "21.78"
$ws.Range("E12").Value = "  -2.70%  "
$ws.Range("D13").Value = "1.911.52"
$ws.Range("E13").Value = "  -4.03%  "
$ws.Range("E14").Value = "  -3.82%  "
$ws.Range("D15").Value = "5.977"
$ws.Range("E15").Value = "  -3.99%  "
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").Value = "89.39"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "0.00001054"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("D19").Value = "0.06585"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").Value = "17.48"
$ws.Range("E20").Value = "  -6.03%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").Value = "5.637"
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("D23").Value = "27.942.42"
$ws.Range("E23").Value = "  -2.98%  "
$ws.Range("D24").Value = "11.06"
$ws.Range("E24").Value = "  -4.10%  "
$ws.Range("D25").Value = "2.302"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").Value = "2.101.56"
$ws.Range("E26").Value = "  -5.47%  "
$ws.Range("D27").Value = "154.14"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "19.88"
$ws.Range("E28").Value = "  -1.50%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "5.702"
$ws.Range("E29").Value = "  -3.72%  "
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "2.106"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("D31").Value = "123.87"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "0.09527"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "0.9588"
$ws.Range("E33").Value = "  -4.50%  "
$ws.Range("D34").Value = "1.476"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("D35").Value = "3.622"
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("D36").Value = "5.462"
$ws.Range("E36").Value = "  -3.75%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "1.256"
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02276"
$ws.Range("E38").Value = "  -2.76%  "
$ws.Range("D39").Value = "0.06100"
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("D40").Value = "8.601"
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("D41").Value = "0.6103"
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "10.73"
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("D44").Value = "0.1895"
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("D45").Value = "1.307"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("D47").Value = "12.75"
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("D48").Value = "1.989"
$ws.Range("E48").Value = "  -3.87%  "
$ws.Range("D49").Value = "3.423"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "0.06894"
$ws.Range("E50").Value = "  +0.83%  "
$ws.Range("D51").Value = "110.32"
$ws.Range("E51").Value = "  -0.68%  "
